$wb = $excel.ActiveWorkbook

# --- RB sheet: "M.Carter" -> "Mi.Carter" (Week 17 roster correction) ---
$wsRB = $wb.Worksheets.Item("RB")
$wsRB.Range("A4").Value = "Mi.Carter"

# --- WR sheet: log Week 17 data - add new player K.Yeboah as row 11 ---
$wsWR = $wb.Worksheets.Item("WR")
$wsWR.Range("A11").Value = "K.Yeboah"
$wsWR.Range("B11:J11").Value = 0

# Make WR the active sheet/tab, with the new selection left on J12
$wsWR.Activate()
$wsWR.Range("J12").Select()
